$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column (H) to the sheet, copying the header formatting
# from the neighboring "sum" header cell (G1) so it matches the existing
# bold/bordered/centered header style, then overwrite the value with "Save".
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Fill in the data values for the new "Save" column.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
